# Update "想去人数" (want-to-go count) figures in column F on both the
# "展览" and "全部类型" sheets, matching the regenerated site data.

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 19
    6  = 1752
    8  = 736
    18 = 4073
    24 = 966
    26 = 26
    28 = 1844
    29 = 47
    33 = 19
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
